$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell C7: italic note "(Pulled Directly from Meeting Minutes)" matching
# the style used for the other italic notes in row 7 (D7/E7).
$ws.Range("C7").Value = "(Pulled Directly from Meeting Minutes)"
$ws.Range("C7").Font.Italic = $true

# Widen column C to match column B's width, mirroring the author's manual
# column resize (the two columns end up the same width).
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Scroll the view right and move the selection, matching the saved
# worksheet view state (topLeftCell B1 / selection E16).
$ws.Range("E16").Select()
$excel.ActiveWindow.ScrollColumn = 2
